$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rewrite the data table (rows 2-33): Model rows moved up, Reference
# rows re-ordered + extended with the final experimental data/statistics ---
$ws.Range("A2").Value = "Glass"
$ws.Range("B2").Value = "Dry"
$ws.Range("C2").Value = "Model"
$ws.Range("D2").Value = 40.299999999999997
$ws.Range("A3").Value = "Glass"
$ws.Range("B3").Value = "Wet"
$ws.Range("C3").Value = "Model"
$ws.Range("D3").Value = 12.4
$ws.Range("A4").Value = "PTFE"
$ws.Range("B4").Value = "Dry"
$ws.Range("C4").Value = "Model"
$ws.Range("D4").Value = 20.3
$ws.Range("A5").Value = "PTFE"
$ws.Range("B5").Value = "Wet"
$ws.Range("C5").Value = "Model"
$ws.Range("D5").Value = 33.85
$ws.Range("A6").Value = "PMMA"
$ws.Range("B6").Value = "Dry"
$ws.Range("C6").Value = "Model"
$ws.Range("D6").Value = 37.99
$ws.Range("A7").Value = "PMMA"
$ws.Range("B7").Value = "Wet"
$ws.Range("C7").Value = "Model"
$ws.Range("D7").Value = 40.78
$ws.Range("A8").Value = "OTS-SAM"
$ws.Range("B8").Value = "Dry"
$ws.Range("C8").Value = "Model"
$ws.Range("D8").Value = 38.97
$ws.Range("A9").Value = "OTS-SAM"
$ws.Range("B9").Value = "Wet"
$ws.Range("C9").Value = "Model"
$ws.Range("D9").Value = 40.06
$ws.Range("A10").Value = "Glass"
$ws.Range("B10").Value = "Dry"
$ws.Range("C10").Value = "Reference"
$ws.Range("D10").Value = 17
$ws.Range("A11").Value = "Glass"
$ws.Range("B11").Value = "Wet"
$ws.Range("C11").Value = "Reference"
$ws.Range("D11").Value = 5.5
$ws.Range("A12").Value = "PMMA"
$ws.Range("B12").Value = "Dry"
$ws.Range("C12").Value = "Reference"
$ws.Range("D12").Value = 27
$ws.Range("A13").Value = "PMMA"
$ws.Range("B13").Value = "Wet"
$ws.Range("C13").Value = "Reference"
$ws.Range("D13").Value = 24
$ws.Range("A14").Value = "OTS-SAM"
$ws.Range("B14").Value = "Dry"
$ws.Range("C14").Value = "Reference"
$ws.Range("D14").Value = 20
$ws.Range("A15").Value = "OTS-SAM"
$ws.Range("B15").Value = "Wet"
$ws.Range("C15").Value = "Reference"
$ws.Range("D15").Value = 17.5
$ws.Range("A16").Value = "PTFE"
$ws.Range("B16").Value = "Dry"
$ws.Range("C16").Value = "Reference"
$ws.Range("D16").Value = 2
$ws.Range("A17").Value = "PTFE"
$ws.Range("B17").Value = "Wet"
$ws.Range("C17").Value = "Reference"
$ws.Range("D17").Value = 7.5
$ws.Range("A18").Value = "Glass"
$ws.Range("B18").Value = "Dry"
$ws.Range("C18").Value = "Reference"
$ws.Range("D18").Value = 21
$ws.Range("A19").Value = "Glass"
$ws.Range("B19").Value = "Wet"
$ws.Range("C19").Value = "Reference"
$ws.Range("D19").Value = 7
$ws.Range("A20").Value = "PMMA"
$ws.Range("B20").Value = "Dry"
$ws.Range("C20").Value = "Reference"
$ws.Range("D20").Value = 28.5
$ws.Range("A21").Value = "PMMA"
$ws.Range("B21").Value = "Wet"
$ws.Range("C21").Value = "Reference"
$ws.Range("D21").Value = 28
$ws.Range("A22").Value = "OTS-SAM"
$ws.Range("B22").Value = "Dry"
$ws.Range("C22").Value = "Reference"
$ws.Range("D22").Value = 21.25
$ws.Range("A23").Value = "OTS-SAM"
$ws.Range("B23").Value = "Wet"
$ws.Range("C23").Value = "Reference"
$ws.Range("D23").Value = 19
$ws.Range("A24").Value = "PTFE"
$ws.Range("B24").Value = "Dry"
$ws.Range("C24").Value = "Reference"
$ws.Range("D24").Value = 2.5
$ws.Range("A25").Value = "PTFE"
$ws.Range("B25").Value = "Wet"
$ws.Range("C25").Value = "Reference"
$ws.Range("D25").Value = 8.5
$ws.Range("A26").Value = "Glass"
$ws.Range("B26").Value = "Dry"
$ws.Range("C26").Value = "Reference"
$ws.Range("D26").Value = 13
$ws.Range("A27").Value = "Glass"
$ws.Range("B27").Value = "Wet"
$ws.Range("C27").Value = "Reference"
$ws.Range("D27").Value = 4
$ws.Range("A28").Value = "PMMA"
$ws.Range("B28").Value = "Dry"
$ws.Range("C28").Value = "Reference"
$ws.Range("D28").Value = 25.5
$ws.Range("A29").Value = "PMMA"
$ws.Range("B29").Value = "Wet"
$ws.Range("C29").Value = "Reference"
$ws.Range("D29").Value = 20
$ws.Range("A30").Value = "OTS-SAM"
$ws.Range("B30").Value = "Dry"
$ws.Range("C30").Value = "Reference"
$ws.Range("D30").Value = 18.75
$ws.Range("A31").Value = "OTS-SAM"
$ws.Range("B31").Value = "Wet"
$ws.Range("C31").Value = "Reference"
$ws.Range("D31").Value = 16
$ws.Range("A32").Value = "PTFE"
$ws.Range("B32").Value = "Dry"
$ws.Range("C32").Value = "Reference"
$ws.Range("D32").Value = 1.5
$ws.Range("A33").Value = "PTFE"
$ws.Range("B33").Value = "Wet"
$ws.Range("C33").Value = "Reference"
$ws.Range("D33").Value = 6.5

# --- Relabeled contact type / surface chemistry => refresh the hidden
#     filter-database defined name scoped to this sheet ---
$fd = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$D`$9")
$fd.Visible = $false

# --- Fixed plot formatting: move the active selection ---
$ws.Range("H8").Select()
